# Apply the edits described by the commit:
#  - "About" sheet: bump the last-updated date in C1 (serial date 45236 -> 45265)
#  - "BGDPbES" sheet: the "biomass" row (row 10) and the "geothermal" row (row 11)
#    are now fully guaranteed-dispatch (value 1) across every year column B:AK,
#    overwriting what used to be 0 / shared formulas referencing column B.
#  - Leave the BGDPbES sheet's selection on cell A10 (where the edit was made).

$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date -----------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45265

# --- BGDPbES sheet: set biomass (row 10) and geothermal (row 11) to 1 --------
$wsB = $wb.Worksheets.Item("BGDPbES")
$wsB.Range("B10:AK10").Value = 1
$wsB.Range("B11:AK11").Value = 1

# Reflect the final cell selection left on the sheet after the edit.
$wsB.Activate()
$wsB.Range("A10").Select()
